$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FBS")

# wind_dir_fg for Marshall @ Louisiana (row 2): NNE -> N
$ws.Range("Q2").Value = "N"

# Odds_n for UNLV @ Boise State (row 3): -110 -> -108
$ws.Range("Z3").Value = -108

# Current for UNLV @ Boise State (row 3): -3.5 -> -4
$ws.Range("AB3").Value = -4

# Move_s for UNLV @ Boise State (row 3): -0.5 -> 0
$ws.Range("AF3").Value = 0

# Timestamp refresh (shared by both data rows)
$ws.Range("AK2").Value = "2024-12-02T10:01:07.532337"
$ws.Range("AK3").Value = "2024-12-02T10:01:07.532337"
